# UI Testing : Personal Detail's Page UI Check Test.
# Adds a new worksheet "T_MyInfoPDUITest" holding the MyInfo/Personal-Details
# field-level UI test data, formats it as a bordered table (Table1), and
# updates the active-sheet/selection state left behind by the edit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The login-test sheet's selection moved to B2 (and it stops being the
# tab-selected sheet once the new sheet becomes active below).
$ws1.Range("B2").Select()

# Insert the new sheet right after T_LoginTest.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "T_MyInfoPDUITest"

$data = @(
    @("Case_ID",  "Field",                  "Status Expected"),
    @("Case_001", "FirstName",              "Enable"),
    @("Case_002", "MiddleName",             "Enable"),
    @("Case_003", "LastName",               "Enable"),
    @("Case_004", "EmployeeId",             "Disable"),
    @("Case_005", "OtherId",                "Enable"),
    @("Case_006", "Driver License Number",  "Disable"),
    @("Case_007", "License Expiry Date",    "Enable"),
    @("Case_008", "Nationality",            "Enable"),
    @("Case_009", "Marital Status",         "Enable"),
    @("Case_010", "Date Of Birth",          "Disable"),
    @("Case_011", "Male Checkbox",          "Enable"),
    @("Case_012", "Female Checkbox",        "Enable")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $data[$i][0]
    $ws2.Cells.Item($row, 2).Value = $data[$i][1]
    $ws2.Cells.Item($row, 3).Value = $data[$i][2]
}

# Turn the Field/Status Expected range into an Excel Table (ListObject) -
# headers are auto-detected from row 1 ("Field" / "Status Expected").
$lo = $ws2.ListObjects.Add(1, $ws2.Range("B1:C13"))

# Thin black border around every cell (reuses the same border style already
# used on T_LoginTest) - set Color before LineStyle so it matches the
# existing style instead of minting a new one.
$rng = $ws2.Range("A1:C13")
$rng.Borders.Color = 0
$rng.Borders.LineStyle = 1

# Size the Field / Status Expected columns to fit their content.
$ws2.Columns("B:C").AutoFit() | Out-Null

# Final selection left on the new sheet.
$ws2.Range("B6").Select()
